$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A39").Value = "2025-04-29 03:14:56"
$ws.Range("B39").Value = 98
